# Rebuild paragraph 1 ("H csif eakcvhisakjwcshi ") so that it:
#   - carries an All-Caps paragraph-mark run property (w:pPr/w:rPr/w:caps)
#   - keeps the original run text as its own run
#   - gets a new run containing just a manual line break
#   - gets a new run containing a manual line break followed by "sumanth github"
#
# We do this with Range.InsertXML, which lets us drop exact OOXML in place
# (so the three <w:r> runs land exactly as required instead of Word's
# Find/Replace or InsertBreak auto-merging/splitting them differently).

$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1).Range

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:caps/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">H csif eakcvhisakjwcshi </w:t>
            </w:r>
            <w:r>
              <w:br/>
            </w:r>
            <w:r>
              <w:br/>
              <w:t>sumanth github</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

# InsertXML inserts the <w:p> above as a *new* paragraph right before the
# original one (the original paragraph, now empty of text, survives as the
# following paragraph). Collapse the two paragraph marks together afterwards
# so only the freshly-built paragraph remains.
$null = $p1.InsertXML($xml)

$firstEnd = $d.Paragraphs(1).Range.End
$secondEnd = $d.Paragraphs(2).Range.End
$null = $d.Range($firstEnd - 1, $secondEnd).Delete()
